# Update the "Lib_1" sheet data: recompute Gap (E) and Time (F) columns,
# and refresh a handful of Z (D) values, reflecting the latest model run
# results for each instance (p1..p57).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 0
$ws.Range("F2").Value = 661
$ws.Range("F3").Value = 563
$ws.Range("F4").Value = 614
$ws.Range("D5").Value = 7168.000000000001
$ws.Range("E5").Value = 0.0000000000009094947017729282
$ws.Range("F5").Value = 1624
$ws.Range("F6").Value = 1223
$ws.Range("F7").Value = 1935
$ws.Range("D8").Value = 4365.999999999998
$ws.Range("E8").Value = 0.000000000001818989403545856
$ws.Range("F8").Value = 4567
$ws.Range("D9").Value = 7926
$ws.Range("E9").Value = 0
$ws.Range("F9").Value = 12981
$ws.Range("F10").Value = 2880
$ws.Range("E11").Value = 0
$ws.Range("F11").Value = 2674
$ws.Range("D12").Value = 3447
$ws.Range("E12").Value = 0.0000000000004547473508864641
$ws.Range("F12").Value = 8166
$ws.Range("F13").Value = 3923
$ws.Range("D14").Value = 3760
$ws.Range("E14").Value = 0
$ws.Range("F14").Value = 5285
$ws.Range("D15").Value = 5965
$ws.Range("E15").Value = 0
$ws.Range("F15").Value = 5831
$ws.Range("D16").Value = 7816
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 5598
$ws.Range("D17").Value = 11543
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 6643
$ws.Range("D18").Value = 9884
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 8780
$ws.Range("D19").Value = 15607
$ws.Range("E19").Value = 0
$ws.Range("F19").Value = 10875
$ws.Range("F20").Value = 11116
$ws.Range("D21").Value = 26561.00000000001
$ws.Range("E21").Value = 0.000000000007275957614183426
$ws.Range("F21").Value = 72196
$ws.Range("D22").Value = 7295
$ws.Range("E22").Value = 0
$ws.Range("F22").Value = 40366
$ws.Range("D23").Value = 3271
$ws.Range("E23").Value = 0
$ws.Range("F23").Value = 5419
$ws.Range("D24").Value = 6036
$ws.Range("E24").Value = 0
$ws.Range("F24").Value = 3418
$ws.Range("D25").Value = 6327
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = 4441
$ws.Range("F26").Value = 2162
$ws.Range("D27").Value = 4448.000000000002
$ws.Range("E27").Value = 0.000000000001818989403545856
$ws.Range("F27").Value = 11238
$ws.Range("D28").Value = 10921
$ws.Range("E28").Value = 0
$ws.Range("F28").Value = 7125
$ws.Range("D29").Value = 11117
$ws.Range("E29").Value = 0
$ws.Range("F29").Value = 15268
$ws.Range("D30").Value = 9832
$ws.Range("E30").Value = 0
$ws.Range("F30").Value = 14677
$ws.Range("D31").Value = 10816
$ws.Range("E31").Value = 0
$ws.Range("F31").Value = 871670
$ws.Range("D32").Value = 4466
$ws.Range("E32").Value = 0
$ws.Range("F32").Value = 17200
$ws.Range("D33").Value = 9880.999999999998
$ws.Range("E33").Value = 0.000000000001818989403545856
$ws.Range("F33").Value = 15855
$ws.Range("D34").Value = 39463
$ws.Range("E34").Value = 0
$ws.Range("F34").Value = 1199075
$ws.Range("D35").Value = 4701
$ws.Range("E35").Value = 0
$ws.Range("F35").Value = 19753
$ws.Range("F36").Value = 2086
$ws.Range("D37").Value = 16781
$ws.Range("E37").Value = 0
$ws.Range("F37").Value = 36640
$ws.Range("D38").Value = 14668
$ws.Range("E38").Value = 0
$ws.Range("F38").Value = 79614
$ws.Range("D39").Value = 47249
$ws.Range("E39").Value = 0
$ws.Range("F39").Value = 20586
$ws.Range("D40").Value = 41007
$ws.Range("E40").Value = 0
$ws.Range("F40").Value = 88409
$ws.Range("D41").Value = 61636
$ws.Range("E41").Value = 3
$ws.Range("F41").Value = 1199401
$ws.Range("E42").Value = 0
$ws.Range("F42").Value = 2998
$ws.Range("D43").Value = 7887
$ws.Range("E43").Value = 0
$ws.Range("F43").Value = 10463
$ws.Range("D44").Value = 5114
$ws.Range("E44").Value = 0
$ws.Range("F44").Value = 20970
$ws.Range("D45").Value = 37303
$ws.Range("E45").Value = 1281
$ws.Range("F45").Value = 1199242
$ws.Range("D46").Value = 17676
$ws.Range("E46").Value = 0
$ws.Range("F46").Value = 4633
$ws.Range("D47").Value = 48701
$ws.Range("E47").Value = 0
$ws.Range("F47").Value = 177261
$ws.Range("D48").Value = 66230
$ws.Range("E48").Value = 0
$ws.Range("F48").Value = 659549
$ws.Range("D49").Value = 58964
$ws.Range("E49").Value = 0
$ws.Range("F49").Value = 65225
$ws.Range("D50").Value = 79659
$ws.Range("E50").Value = 45
$ws.Range("F50").Value = 1199777
$ws.Range("D51").Value = 5937
$ws.Range("E51").Value = 0
$ws.Range("F51").Value = 29383
$ws.Range("D52").Value = 9060
$ws.Range("E52").Value = 0
$ws.Range("F52").Value = 224940
$ws.Range("D53").Value = 34652
$ws.Range("E53").Value = 0
$ws.Range("F53").Value = 560694
$ws.Range("D54").Value = 30038
$ws.Range("E54").Value = 0
$ws.Range("F54").Value = 33280
$ws.Range("D55").Value = 43853
$ws.Range("E55").Value = 0
$ws.Range("F55").Value = 10129
$ws.Range("D56").Value = 69753
$ws.Range("E56").Value = 143
$ws.Range("F56").Value = 1199989
$ws.Range("D57").Value = 64478
$ws.Range("E57").Value = 4
$ws.Range("F57").Value = 1155
$ws.Range("F58").Value = 611

Write-Host "Updated cells with latest instance data."
